# Update the publication Date on the Metadata sheet
$wb = $excel.ActiveWorkbook
$metaWs = $wb.Worksheets.Item("Metadata")
$metaWs.Range("B8").Value = "2024-03-22T16:25:12+00:00"

# On the Elements sheet, the two "Mapping" columns (AK = column 37 and
# AL = column 38) have been swapped: "Mapping: Spécification métier vers
# l'extension ROR LocationSupportedCapacity" now comes before
# "Mapping: RIM Mapping". Swap the whole columns (header + all data rows)
# by cutting column AL and inserting it in front of column AK.
$ws = $wb.Worksheets.Item("Elements")
$ws.Columns.Item(38).Cut()
$ws.Columns.Item(37).Insert()

# The Insert() operation resets the column widths to the default, so
# restore the (swapped) widths that used to belong to each column.
# Former AL width (83.625) now belongs to AK (col 37); former AK width
# (24.98046875) now belongs to AL (col 38).
$ws.Columns.Item(37).ColumnWidth = 82.83333333333333
$ws.Columns.Item(38).ColumnWidth = 24.166666666666668

# Columns C, D, AE, AF, AG are hidden helper columns in the original
# workbook and are unaffected by this change; re-assert their hidden
# state so the round trip through the engine keeps them hidden.
$ws.Columns.Item(3).Hidden = $true
$ws.Columns.Item(4).Hidden = $true
$ws.Columns.Item(31).Hidden = $true
$ws.Columns.Item(32).Hidden = $true
$ws.Columns.Item(33).Hidden = $true
